$d = $word.ActiveDocument

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/><w:lang w:val="en-US"/></w:rPr><w:t>Thought Process and Approach for Google AdSense Checker</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="767171" w:themeColor="background2" w:themeShade="80"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Initial Considerations</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>The first thing I thought about was how to spot the signs of Google AdSense on a website. Google AdSense leaves multiple indications in the site's source code, which can be used to identify its presence. These indicators include:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>- `adsbygoogle.js`</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>- `googlesyndication.com`</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>- Specific `&lt;div&gt;` elements with class names related to ads</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Identifying Patterns and Indications</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>After understanding how to spot these indicators, I planned to create a scraper function to detect and collect the data needed to classify a site as having Google AdSense.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Challenges and Dynamic Detection</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>I realized that scraping the HTML alone might not always work because a site can support Google AdSense without showing ads all the time. To reliably detect ads, I needed a way to do it dynamically and in real time.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Solution with Selenium</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>That's when I came across Selenium, a tool for web testing that was perfect for this. Selenium allows for real-time interaction with web pages, making it possible to detect dynamically loaded ads.</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>Using iframe detection works for many Google AdSense ads since they are often served within iframe elements.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Developing the Code</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>Here's a breakdown of how I developed the code:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>1. **Preparation: I started by setting up Selenium WebDriver and pandas for handling input and output in CSV format. The CSV files included one for the list of websites and another for the results.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>2. **Data Collection: I used an API to fetch a bulk of websites and append them to the websites file.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>3. Ad Detection: The main part of the code involves detecting ads. I created the `check_adsense_ads` function to initialize the Edge WebDriver, navigate to a URL, and wait up to 30 seconds for an iframe containing Google ads to load. If such an iframe is found, the function returns True, indicating that AdSense ads are present.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>4. Running the Script: The `main` function reads URLs from the input CSV file, checks each URL for ads, and writes the results to an output CSV file. By using `WebDriverWait` with `EC.frame_to_be_available_and_switch_to_it`, the script can wait for dynamically loaded ads, making it robust for websites where ads load after the main content.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>5. Execution: The script uses the standard Python idiom (`if __name__ == "__main__":`) to ensure the `main` function runs only when the script is executed directly, not when imported as a module.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="en-US"/></w:rPr><w:t>This setup ensures that websites are accurately checked for AdSense ads, with results systematically logged in a CSV file for further analysis.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($xml)

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
